$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename existing "_car" headers to "_icev" ---
$ws.Range("E1").Value = "avg_speed_icev"
$ws.Range("F1").Value = "energy_per_mile_icev"
$ws.Range("G1").Value = "avg_veh_num_icev"

# --- Step 2: insert 9 new columns after G (for ev/cav/caev vehicle stats) ---
$ws.Range("H1:P1").EntireColumn.Insert()

$ws.Range("H1").Value = "avg_speed_ev"
$ws.Range("I1").Value = "energy_per_mile_ev"
$ws.Range("J1").Value = "avg_veh_num_ev"
$ws.Range("K1").Value = "avg_speed_cav"
$ws.Range("L1").Value = "energy_per_mile_cav"
$ws.Range("M1").Value = "avg_veh_num_cav"
$ws.Range("N1").Value = "avg_speed_caev"
$ws.Range("O1").Value = "energy_per_mile_caev"
$ws.Range("P1").Value = "avg_veh_num_caev"

# --- Step 3: insert 2 new columns before "M_reward_interpolation" ---
# after the first insertion, old AN (M_reward_interpolation) moved to AW
$ws.Range("AW1:AX1").EntireColumn.Insert()

$ws.Range("AW1").Value = "En_phase_mode"
$ws.Range("AX1").Value = "En_vehicle_type_distribution"
